$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "spherical fog is back for 1.21.6"

$ws.Range("A19").Select()
